# "removed console from slides"
# Removes the "Console" hexagon TextBox (id=53, "TextBox 33") and the elbow
# connector (id=64, "Connector: Elbow 38") that links it to the "Receiver
# service" box on slide 1.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

foreach ($shape in $s.Shapes) {
    if ($shape.Id -eq 64 -or $shape.Id -eq 53) {
        $shape.Delete()
    }
}
